# Add a new "user person/organisation category" criterion row to the
# criterion_property sheet. This inserts one new row at row 275 (pushing
# the existing rows 275-293 down to 276-294) and fills in the new row's
# cells for the USER_DB module's "user.identity.category.id" property.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("criterion_property")

# Insert a brand-new row above the current row 275 ("user.authMethod"),
# shifting it (and everything below) down by one.
$ws.Rows("275:275").Insert()

# Match the row height/customHeight formatting used by every other data
# row on this sheet.
$ws.Rows("275:275").RowHeight = 16.5

# Populate the new row: USER_DB / user.identity.category.id / LONG /
# StaffCategory / getAllStaffCategories / getName / getId / ... / EQ, NE
$ws.Range("A275").Value = "USER_DB"
$ws.Range("B275").Value = "user.identity.category.id"
$ws.Range("C275").Value = "LONG"
$ws.Range("D275").Value = "StaffCategory"
$ws.Range("E275").Value = "getAllStaffCategories"
$ws.Range("F275").Value = "getName"
$ws.Range("G275").Value = "getId"
$ws.Range("K275").Value = "user.identity.category.id"
$ws.Range("L275").Value = "EQ, NE"

# Move the active selection to the newly added row, as in the source
# workbook.
$ws.Range("A275").Select() | Out-Null
